$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (want-to-go count) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 250
$wsExpo.Range("F3").Value = 559
$wsExpo.Range("F5").Value = 268
$wsExpo.Range("F6").Value = 1063
$wsExpo.Range("F7").Value = 1395
$wsExpo.Range("F10").Value = 732
$wsExpo.Range("F12").Value = 125
$wsExpo.Range("F13").Value = 109
$wsExpo.Range("F14").Value = 400
$wsExpo.Range("F15").Value = 1274
$wsExpo.Range("F16").Value = 94
$wsExpo.Range("F17").Value = 76
$wsExpo.Range("F18").Value = 260
$wsExpo.Range("F19").Value = 5215
$wsExpo.Range("F20").Value = 630
$wsExpo.Range("F21").Value = 25
$wsExpo.Range("F22").Value = 186
$wsExpo.Range("F23").Value = 5483
$wsExpo.Range("F24").Value = 47
$wsExpo.Range("F25").Value = 110
$wsExpo.Range("F26").Value = 82
$wsExpo.Range("F28").Value = 13905
$wsExpo.Range("F29").Value = 1404
$wsExpo.Range("F31").Value = 84
$wsExpo.Range("F33").Value = 396
$wsExpo.Range("F34").Value = 547
$wsExpo.Range("F35").Value = 4139
$wsExpo.Range("F36").Value = 88
$wsExpo.Range("F37").Value = 347

# Sheet "演出" (Performances) - update "最低票价" (minimum ticket price) value
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("G4").Value = 210

# Sheet "全部类型" (All types) - update "想去人数" and "最低票价" values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 250
$wsAll.Range("F3").Value = 559
$wsAll.Range("F5").Value = 268
$wsAll.Range("F6").Value = 1063
$wsAll.Range("F7").Value = 1395
$wsAll.Range("F10").Value = 732
$wsAll.Range("F12").Value = 125
$wsAll.Range("F13").Value = 109
$wsAll.Range("F14").Value = 400
$wsAll.Range("F15").Value = 1274
$wsAll.Range("F16").Value = 94
$wsAll.Range("F17").Value = 76
$wsAll.Range("F18").Value = 260
$wsAll.Range("F20").Value = 5215
$wsAll.Range("F21").Value = 630
$wsAll.Range("F23").Value = 25
$wsAll.Range("F24").Value = 186
$wsAll.Range("G25").Value = 210
$wsAll.Range("F26").Value = 5483
$wsAll.Range("F27").Value = 47
$wsAll.Range("F28").Value = 110
$wsAll.Range("F29").Value = 82
$wsAll.Range("F31").Value = 13906
$wsAll.Range("F32").Value = 1404
$wsAll.Range("F34").Value = 84
$wsAll.Range("F36").Value = 396
$wsAll.Range("F37").Value = 547
$wsAll.Range("F38").Value = 4139
$wsAll.Range("F39").Value = 88
$wsAll.Range("F40").Value = 347
